# autosave.xlsx edit — add calculated energy/CO2 fields to the
# "Autosave Fields" sheet, and switch the active tab from "Examples" to
# "Autosave Fields".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autosave Fields")

# --- Insert 9 new rows before the NOTES row (old row 70), pushing
#     everything below it down by 9 rows (old 70 -> new 79, old 88 -> new 97).
#     Excel copies the row-above formatting on insert, which reproduces the
#     existing s="1" (A) / s="5" (C) styling for the new rows automatically.
$ws.Rows("70:78").Insert()

# --- Column A: new internal field names (burner labels) ---
$ws.Range("A70").Value = "~btubatch"
$ws.Range("A71").Value = "~co2batch"
$ws.Range("A72").Value = "~btupreheat"
$ws.Range("A73").Value = "~co2preheat"
$ws.Range("A74").Value = "~btubbp"
$ws.Range("A75").Value = "~co2bbp"
$ws.Range("A76").Value = "~bturoast"
$ws.Range("A77").Value = "~co2roast"
$ws.Range("A78").Value = "~co2pergreenkg"

# --- Column B: descriptions (note: insertion order mirrors the source
#     workbook's shared-string table ordering) ---
$ws.Range("B71").Value = "From the Profile Energy Use - CO2 produced by the batch in g"
$ws.Range("B72").Value = "From the Profile Energy Use - Energy used during preheat in BTU"
$ws.Range("B70").Value = "From the Profile Energy Use - Total energy used by the batch in BTU"
$ws.Range("B76").Value = "From the Profile Energy Use - Energy used from CHARGE to DROP in BTU"
$ws.Range("B74").Value = "From the Profile Energy Use - Energy used during Between Batch Protocol in BTU"
$ws.Range("B77").Value = "From the Profile Energy Use - CO2 produced from CHARGE to DROP in g"
$ws.Range("B75").Value = "From the Profile Energy Use - CO2 produced during Between Batch Protocol in g"
$ws.Range("B73").Value = "From the Profile Energy Use - CO2 produced during preheat in g"
$ws.Range("B78").Value = "From the Profile Energy Use - CO2 produced per kg of green beans in g"

# B column cells wrap like the other long description cells in this sheet.
$ws.Range("B70:B78").WrapText = $true

# --- Column C: example values ---
$ws.Range("C70").Value = 8943.2000000000007
$ws.Range("C71").Value = 923.3
$ws.Range("C72").Value = 2538.8000000000002
$ws.Range("C73").Value = 443.9
$ws.Range("C74").Value = 1019.7
$ws.Range("C75").Value = 254.1
$ws.Range("C76").Value = 7843.2
$ws.Range("C77").Value = 873.9
$ws.Range("C78").Value = 354.3

# --- Switch the active/selected sheet & view from "Examples" to
#     "Autosave Fields", matching the new workbook activeTab + sheetView. ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 62
$ws.Range("C77").Select()

# --- Best-effort: turn on iterative calculation with a small max change,
#     mirroring the new calcPr iterateDelta="1E-4". ---
$excel.Iteration = $true
$excel.MaxChange = 0.0001
